$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix unit strings in column C (trailing space / caret notation) ---
$ws.Range("C2").Value  = "dt/ha"
$ws.Range("C3").Value  = "g/m^2"
$ws.Range("C4").Value  = "#/m^2"
$ws.Range("C21").Value = "g/m^2"
$ws.Range("C22").Value = "dt/ha"
$ws.Range("C24").Value = "dt/ha"
$ws.Range("C25").Value = "#/m^2"
$ws.Range("C26").Value = "dt/ha"

# --- New column F: trait abbreviation ---
$ws.Range("F1").Value  = "abbrev"
$ws.Range("F2").Value  = "GY"
$ws.Range("F4").Value  = "SN"
$ws.Range("F5").Value  = "GpS"
$ws.Range("F6").Value  = "GN"
$ws.Range("F8").Value  = "TKW"
$ws.Range("F22").Value = "Shoot"
$ws.Range("F23").Value = "HI"
$ws.Range("F24").Value = "Straw"
$ws.Range("F25").Value = "GN"
$ws.Range("F26").Value = "GP"

# --- Column C width & current selection, matching the saved view state ---
$ws.Range("C1").ColumnWidth = 21.6
$ws.Range("C4").Select()
